$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the Date value ---
$metaSheet = $wb.Worksheets.Item("Metadata")
$metaSheet.Range("B8").Value = "2026-01-23T08:28:04+00:00"

# --- Mapping Table 1 sheet: body-section mapping updates ---
$mapSheet = $wb.Worksheets.Item("Mapping Table 1")

# Row 8: FRCDADirectiveAnticipee.value -> FRCDADirectiveAnticipee.valueBoolean
$mapSheet.Range("A8").Value = "FRCDADirectiveAnticipee.valueBoolean"

# Insert a new row before the current row 12 (entryRelationship.observationMedia.value row)
# for the new FRCDADirectiveAnticipee.entryRelationship.observationMedia.id mapping.
$mapSheet.Rows.Item(12).Insert()

$mapSheet.Range("A12").Value = "FRCDADirectiveAnticipee.entryRelationship.observationMedia.id"
$mapSheet.Range("C12").Value = "equivalent"
$mapSheet.Range("D12").Value = "FRAdvanceDirectiveDocument.sourceAttachment.id"

# Copy the style of the neighboring data row onto the newly inserted row
$mapSheet.Range("A11:E11").Copy()
$mapSheet.Range("A12:E12").PasteSpecial(-4122) # xlPasteFormats
